# Add three new market test-data sheets (Netherlands, Austria, Denmark)
# right after the existing "Greece" sheet, each built from the "Greece"
# sheet template (same layout/styles/merged cells), with the market name
# (B2) and user-story ticket id (B4) swapped in for each new country.

$wb = $excel.ActiveWorkbook
$greece = $wb.Worksheets.Item("Greece")

# --- Netherlands -----------------------------------------------------
$greece.Copy($null, $greece)
$netherlands = $wb.Worksheets.Item($greece.Index + 1)
$netherlands.Name = "Netherlands"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Range("B4").Value = "NGC-3144/T2192"

# --- Austria -----------------------------------------------------------
$greece.Copy($null, $netherlands)
$austria = $wb.Worksheets.Item($netherlands.Index + 1)
$austria.Name = "Austria"
$austria.Range("B2").Value = "Austria Market"
$austria.Range("B4").Value = "NGC-3817/T2299"

# --- Denmark -----------------------------------------------------------
$greece.Copy($null, $austria)
$denmark = $wb.Worksheets.Item($austria.Index + 1)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-2913/T2779"

# Restore each sheet's "last selected cell" to B4 (where the user-story
# id was typed), then finish with Denmark as the active tab/sheet.
$netherlands.Activate()
[void]$netherlands.Range("B4").Select()

$austria.Activate()
[void]$austria.Range("B4").Select()

$denmark.Activate()
[void]$denmark.Range("B4").Select()

# Scroll the sheet-tab strip so the first visible tab is index 3
# (best-effort - cosmetic only).
$excel.ActiveWindow.ScrollWorkbookTabs([System.Reflection.Missing]::Value, 3)
